$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.3744571182745915
$ws.Range("C2").Value = 0.06005086765733836
$ws.Range("E2").Value = 0.5865357784368257
$ws.Range("F2").Value = 2.449169946602794
$ws.Range("G2").Value = 0.7060045748315815
$ws.Range("H2").Value = 0.8060969680762184
$ws.Range("I2").Value = 0.6189324368316598
$ws.Range("J2").Value = 0.05891706508964489
$ws.Range("K2").Value = 0.4070731114357784
$ws.Range("N2").Value = 1.464706026856127
$ws.Range("B3").Value = 0.336929604762048
$ws.Range("C3").Value = 0.0528572948692414
$ws.Range("E3").Value = 0.5621997698815164
$ws.Range("F3").Value = 2.416051399948927
$ws.Range("G3").Value = 0.70818315275271
$ws.Range("H3").Value = 0.8116852256986249
$ws.Range("I3").Value = 0.6245805938032412
$ws.Range("J3").Value = 0.0592587495646768
$ws.Range("K3").Value = 0.3642351848663168
$ws.Range("N3").Value = 1.482308678207872
$ws.Range("B4").Value = 0.3139622367207835
$ws.Range("C4").Value = 0.04844568282838679
$ws.Range("E4").Value = 0.5475473916211939
$ws.Range("F4").Value = 2.397185469945342
$ws.Range("G4").Value = 0.7100314026908947
$ws.Range("H4").Value = 0.8155048964673384
$ws.Range("I4").Value = 0.6284047086951574
$ws.Range("J4").Value = 0.05950553241511614
$ws.Range("K4").Value = 0.3380054391976159
$ws.Range("N4").Value = 1.493666059897231
$ws.Range("B5").Value = 0.30462197637911
$ws.Range("C5").Value = 0.04664922509469704
$ws.Range("E5").Value = 0.5416492343913291
$ws.Range("F5").Value = 2.389866223045118
$ws.Range("G5").Value = 0.7109126913038821
$ws.Range("H5").Value = 0.8171590901989347
$ws.Range("I5").Value = 0.6300525135975761
$ws.Range("J5").Value = 0.05961538257117382
$ws.Range("K5").Value = 0.3273352313387647
$ws.Range("N5").Value = 1.498432346281201
$ws.Range("B6").Value = 0.303072198958688
$ws.Range("C6").Value = 0.04635100334597553
$ws.Range("E6").Value = 0.5406742436257161
$ws.Range("F6").Value = 2.38867312278812
$ws.Range("G6").Value = 0.7110667580290055
$ws.Range("H6").Value = 0.817439664466427
$ws.Range("I6").Value = 0.6303315297212642
$ws.Range("J6").Value = 0.05963418336187587
$ws.Range("K6").Value = 0.3255645851024553
$ws.Range("N6").Value = 1.499232121343374
$ws.Range("B7").Value = 0.3138361927867379
$ws.Range("C7").Value = 0.04842144988639063
$ws.Range("E7").Value = 0.5474675523845463
$ws.Range("F7").Value = 2.39708526769752
$ws.Range("G7").Value = 0.7100427697127927
$ws.Range("H7").Value = 0.8155268101818294
$ws.Range("I7").Value = 0.628426569539144
$ws.Range("J7").Value = 0.05950697632327007
$ws.Range("K7").Value = 0.3378614612375088
$ws.Range("N7").Value = 1.493729780780715
$ws.Range("B8").Value = 0.3615023303984799
$ws.Range("C8").Value = 0.05756941825995909
$ws.Range("E8").Value = 0.5780844428410035
$ws.Range("F8").Value = 2.437445392468263
$ws.Range("G8").Value = 0.7066496035451593
$ws.Range("H8").Value = 0.8079431537292976
$ws.Range("I8").Value = 0.6208059455985691
$ws.Range("J8").Value = 0.05902719112017607
$ws.Range("K8").Value = 0.3922876219509419
$ws.Range("N8").Value = 1.470661382496023
$ws.Range("B9").Value = 0.455557460649004
$ws.Range("C9").Value = 0.07555199534149892
$ws.Range("E9").Value = 0.6404347368179515
$ws.Range("F9").Value = 2.528282767420762
$ws.Range("G9").Value = 0.7040604629905545
$ws.Range("H9").Value = 0.7961556993468974
$ws.Range("I9").Value = 0.6086923948930689
$ws.Range("J9").Value = 0.05838060134324152
$ws.Range("K9").Value = 0.4995884316618628
$ws.Range("N9").Value = 1.429784666233606
$ws.Range("B10").Value = 0.5250075405654684
$ws.Range("C10").Value = 0.08879384852636463
$ws.Range("E10").Value = 0.6876706627877525
$ws.Range("F10").Value = 2.602208537893631
$ws.Range("G10").Value = 0.7046561492403924
$ws.Range("H10").Value = 0.7893785610372248
$ws.Range("I10").Value = 0.6015249239741856
$ws.Range("J10").Value = 0.05808610094181077
$ws.Range("K10").Value = 0.5787690205867477
$ws.Range("N10").Value = 1.402412539191405
$ws.Range("B11").Value = 0.5566766685415416
$ws.Range("C11").Value = 0.09482528959040337
$ws.Range("E11").Value = 0.7094737653319498
$ws.Range("F11").Value = 2.637413765676911
$ws.Range("G11").Value = 0.7054739208137875
$ws.Range("H11").Value = 0.7867050623620031
$ws.Range("I11").Value = 0.5986419246156025
$ws.Range("J11").Value = 0.05799158504154889
$ws.Range("K11").Value = 0.6148657943774936
$ws.Range("N11").Value = 1.390538258310817
$ws.Range("B12").Value = 0.5686796144545951
$ws.Range("C12").Value = 0.09711038000648387
$ws.Range("E12").Value = 0.7177756044932835
$ws.Range("F12").Value = 2.65097260581058
$ws.Range("G12").Value = 0.7058625549753259
$ws.Range("H12").Value = 0.7857516144120353
$ws.Range("I12").Value = 0.5976046200893883
$ws.Range("J12").Value = 0.05796148797446321
$ws.Range("K12").Value = 0.6285456238183826
$ws.Range("N12").Value = 1.386124890029341
$ws.Range("B13").Value = 0.5660941028900766
$ws.Range("C13").Value = 0.09661819510489522
$ws.Range("E13").Value = 0.7159856299133622
$ws.Range("F13").Value = 2.648042341906148
$ws.Range("G13").Value = 0.7057753384772667
$ws.Range("H13").Value = 0.7859543336514179
$ws.Range("I13").Value = 0.5978255994166446
$ws.Range("J13").Value = 0.05796771634162567
$ws.Range("K13").Value = 0.6255989542015641
$ws.Range("N13").Value = 1.387071687002137
$ws.Range("B14").Value = 0.5576639492014976
$ws.Range("C14").Value = 0.09501326293568013
$ws.Range("E14").Value = 0.7101558503538854
$ws.Range("F14").Value = 2.638524698143868
$ws.Range("G14").Value = 0.705504309551273
$ws.Range("H14").Value = 0.7866254398218615
$ws.Range("I14").Value = 0.5985554934014914
$ws.Range("J14").Value = 0.05798899470145713
$ws.Range("K14").Value = 0.6159910283025454
$ws.Range("N14").Value = 1.390173500357255
$ws.Range("B15").Value = 0.5525015971354321
$ws.Range("C15").Value = 0.09403034120748544
$ws.Range("E15").Value = 0.7065908689288847
$ws.Range("F15").Value = 2.632724501208884
$ws.Range("G15").Value = 0.7053485895966105
$ws.Range("H15").Value = 0.7870441905570118
$ws.Range("I15").Value = 0.5990096663617308
$ws.Range("J15").Value = 0.05800277046041913
$ws.Range("K15").Value = 0.6101072884837038
$ws.Range("N15").Value = 1.392084284937646
$ws.Range("B16").Value = 0.5229393836719396
$ws.Range("C16").Value = 0.08839983488152825
$ws.Range("E16").Value = 0.6862521346502746
$ws.Range("F16").Value = 2.599939571368225
$ws.Range("G16").Value = 0.7046137372981178
$ws.Range("H16").Value = 0.7895615256057624
$ws.Range("I16").Value = 0.6017209431875266
$ws.Range("J16").Value = 0.05809307345138137
$ws.Range("K16").Value = 0.5764115356013519
$ws.Range("N16").Value = 1.403200181859646
$ws.Range("B17").Value = 0.5048230982103803
$ws.Range("C17").Value = 0.08494767202932962
$ws.Range("E17").Value = 0.6738558239980676
$ws.Range("F17").Value = 2.580231270334679
$ws.Range("G17").Value = 0.7043032127614026
$ws.Range("H17").Value = 0.7912107490551819
$ws.Range("I17").Value = 0.6034810133155801
$ws.Range("J17").Value = 0.05815859112206567
$ws.Range("K17").Value = 0.5557598040867049
$ws.Range("N17").Value = 1.410167410948727
$ws.Range("B18").Value = 0.4944102605901151
$ws.Range("C18").Value = 0.0829627990031554
$ws.Range("E18").Value = 0.6667554572240277
$ws.Range("F18").Value = 2.569043874065585
$ws.Range("G18").Value = 0.7041760537695581
$ws.Range("H18").Value = 0.7921978744572726
$ws.Range("I18").Value = 0.6045288846785226
$ws.Range("J18").Value = 0.05819998749100108
$ws.Range("K18").Value = 0.543888753267936
$ws.Range("N18").Value = 1.414229119951059
$ws.Range("B19").Value = 0.4908858960160103
$ws.Range("C19").Value = 0.08229087830626725
$ws.Range("E19").Value = 0.6643564853327319
$ws.Range("F19").Value = 2.565281461274708
$ws.Range("G19").Value = 0.7041418257531404
$ws.Range("H19").Value = 0.7925387146282503
$ws.Range("I19").Value = 0.6048897732819682
$ws.Range("J19").Value = 0.05821464056057835
$ws.Range("K19").Value = 0.5398706795598684
$ws.Range("N19").Value = 1.415613674434296
$ws.Range("B20").Value = 0.506750868388167
$ws.Range("C20").Value = 0.08531508607131855
$ws.Range("E20").Value = 0.6751723613497091
$ws.Range("F20").Value = 2.582313898974434
$ws.Range("G20").Value = 0.7043309416895767
$ws.Range("H20").Value = 0.7910311976064435
$ws.Range("I20").Value = 0.6032899730825996
$ws.Range("J20").Value = 0.05815123229216113
$ws.Range("K20").Value = 0.5579574651871155
$ws.Range("N20").Value = 1.409420112280125
$ws.Range("B21").Value = 0.5601398068522201
$ws.Range("C21").Value = 0.09548464005757751
$ws.Range("E21").Value = 0.7118669625040042
$ws.Range("F21").Value = 2.641314083628259
$ws.Range("G21").Value = 0.7055817716487809
$ws.Range("H21").Value = 0.7864267192866663
$ws.Range("I21").Value = 0.5983396276690556
$ws.Range("J21").Value = 0.05798259003123718
$ws.Range("K21").Value = 0.6188128194814055
$ws.Range("N21").Value = 1.389260164482957
$ws.Range("B22").Value = 0.595093789477886
$ws.Range("C22").Value = 0.1021375347950197
$ws.Range("E22").Value = 0.7361141528174926
$ws.Range("F22").Value = 2.681199841383233
$ws.Range("G22").Value = 0.7068596669698621
$ws.Range("H22").Value = 0.783761053229199
$ws.Range("I22").Value = 0.5954215845544653
$ws.Range("J22").Value = 0.05790557154820775
$ws.Range("K22").Value = 0.6586479055969789
$ws.Range("N22").Value = 1.376569280449687
$ws.Range("B23").Value = 0.5764327159953382
$ws.Range("C23").Value = 0.09858615996375875
$ws.Range("E23").Value = 0.7231486609806836
$ws.Range("F23").Value = 2.65979050643719
$ws.Range("G23").Value = 0.7061353924232492
$ws.Range("H23").Value = 0.7851523044708699
$ws.Range("I23").Value = 0.5969499208800144
$ws.Range("J23").Value = 0.05794363319458284
$ws.Range("K23").Value = 0.6373815504166771
$ws.Range("N23").Value = 1.383298240459073
$ws.Range("B24").Value = 0.5058793158109722
$ws.Range("C24").Value = 0.08514897874113103
$ws.Range("E24").Value = 0.67457707247695
$ws.Range("F24").Value = 2.581371896650978
$ws.Range("G24").Value = 0.7043182454610957
$ws.Range("H24").Value = 0.7911122514432662
$ws.Range("I24").Value = 0.6033762302919499
$ws.Range("J24").Value = 0.05815454760433525
$ws.Range("K24").Value = 0.556963896647062
$ws.Range("N24").Value = 1.409757791086501
$ws.Range("B25").Value = 0.430051481700616
$ws.Range("C25").Value = 0.07068221371231687
$ws.Range("E25").Value = 0.6233180538699514
$ws.Range("F25").Value = 2.502451021687989
$ws.Range("G25").Value = 0.7043236729186049
$ws.Range("H25").Value = 0.7990140445613747
$ws.Range("I25").Value = 0.6116657080157246
$ws.Range("J25").Value = 0.05852389871903085
$ws.Range("K25").Value = 0.4704997303997516
$ws.Range("N25").Value = 1.440375923749018
